# Apply the cryptos-list refresh edits (value updates for rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.892.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.01%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.872.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.88%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7411"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.05%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9992"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3148"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.77%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07213"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.66%  "

# Row 10
$ws.Range("E10").Value = "  -3.84%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08303"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.46%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7500"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.61%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.901.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.68%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.377"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.33%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.58%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.107"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.45%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.889.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.62%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "246.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.25%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.66%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007833"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.36%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9989"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.131.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.12%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1541"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.72%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.281"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.27%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.96%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.47%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.013"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.499"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.12%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.585"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.09%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.535"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.01%  "

# Row 33
$ws.Range("E33").Value = "  +3.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05324"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.22%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.238"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.35%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7492"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.02%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.000"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.697"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.03%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01963"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.83%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.752"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.99%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4516"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.20%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.112.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.19%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.127"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.07%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.58%  "

# Row 45
$ws.Range("E45").Value = "  +1.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.36%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.06%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.862"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.16%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.612"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.41%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.513"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.56%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.029.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.68%  "
